# Burndown chart and sprint log grammar fixes
# Fix the "Identifie" -> "Identify" typo in the sprint task descriptions,
# and move the saved selection to the last-edited cell (C8), matching the
# author's final view state when the workbook was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burndown Chart")

$ws.Range("C7").Value = "Identify GoF design patterns"
$ws.Range("C8").Value = "Identify code smells"
$ws.Range("C9").Value = "Review each others work"
$ws.Range("C10").Value = "Make report"

# Leave the selection / saved view on C8 (no frozen/scrolled topLeftCell),
# as in the committed workbook.
$ws.Range("C8").Select()
